# MLB Final Project w Bugs
# - Clear the AVERAGE() summary formulas (and their cached results) that
#   lived in row 32 (C32:F32), turning them back into blank (but still
#   styled) cells.
# - Reset the sheet view: select A1 / scroll back to the top-left instead
#   of the stale "topLeftCell=A14 / selection=E2" view state left over
#   from the previous editing session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the average formulas from the bottom "totals" row, leaving the
# cells empty but keeping their existing number/alignment formatting.
$ws.Range("C32:F32").ClearContents()

# Reset the view back to the top of the sheet with A1 selected.
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "cleared C32:F32 formulas and reset sheet view"
